# Insert a new weekly price record at row 171 for
# "Feria Lagunitas de Puerto Montt - Pepino ensalada", pushing the
# existing rows 171..271 down to 172..272 (dimension grows from
# A1:R271 to A1:R272).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 171 downward (this also grows the used range / dimension).
$ws.Rows("171").Insert()

# Seed the newly-inserted (blank) row 171 with the same row template as
# the row that is now directly below it (former row 171, now row 172),
# then overwrite the columns that hold the new record's own data.
$ws.Range("A172:R172").Copy()
$ws.Range("A171:R171").PasteSpecial()

$ws.Cells.Item(171, 4).Value = 44719   # D171 Fecha
$ws.Cells.Item(171, 10).Value = 400    # J171 Volumen
$ws.Cells.Item(171, 11).Value = 23000  # K171 Precio minimo
$ws.Cells.Item(171, 12).Value = 23500  # L171 Precio maximo
$ws.Cells.Item(171, 13).Value = 23250  # M171 Precio promedio ponderado
$ws.Cells.Item(171, 16).Value = 388    # P171 Precio $/Kg
